$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Expand the "Declaration" paragraph wording (Self Help Community Housing
#    Association / Data Protection Act -> The C.H.E.E.S.E. Project / EU GDPR).
# ---------------------------------------------------------------------------
$oldDeclaration = "I declare that the information given in this application is, to the best of my knowledge, complete and correct and that it may be used for purposes registered by Self Help Community Housing Association under the Data Protection Act.  I understand that if, after appointment, any information is found to be inaccurate this may lead to disciplinary action or dismissal without notice."

$newDeclaration = "I declare that the information given in this application is, to the best of my knowledge, complete and correct and that it may be used for purposes outlined by The C.H.E.E.S.E. Project for its staff and members under the EU GDPR (full details here: https://cheeseproject.co.uk/privacy-notice-staff-and-members).  I understand that if, after appointment, any information is found to be inaccurate this may lead to disciplinary action or dismissal without notice."

$declRange = $d.Content
$declRange.Find.Execute($oldDeclaration, $true, $false, $false, $false, $false, $true, 1, $false, $newDeclaration, 2)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from after the "Cold Homes Energy Efficiency
#    Survey Experts" subtitle to the end of the (now rewritten) Declaration
#    paragraph. Bookmarks.Add re-uses the reserved name, so the stale
#    bookmark pair that used to sit after the subtitle is dropped
#    automatically once the new one is created.
#
#    A zero-length range placed exactly at a paragraph's trailing text
#    boundary gets normalised to span the whole paragraph, so a short-lived
#    marker run is inserted first, the bookmark is anchored around it, and
#    the marker text is then deleted -- leaving a clean, collapsed bookmark
#    immediately after the final run and before the paragraph mark.
# ---------------------------------------------------------------------------
$marker = "@@GoBackMarker@@"

$tail = $d.Content
$tail.Find.Execute("dismissal without notice.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Collapse(0)
$tail.InsertAfter($marker)

$markerRange = $d.Content
$markerRange.Find.Execute("notice." + $marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange.Collapse(0)
$markerRange.MoveStart(1, -1 * $marker.Length)
$d.Bookmarks.Add("_GoBack", $markerRange)

$cleanup = $d.Content
$cleanup.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cleanup.Delete()
